$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the worker debt table (rows 16-43): previous account statement
# entries removed, new ones added per updated database (NIT-9006817920).
$rows = @(
    ,@(16, '73131545', 'LUIS E ROMERO CABARCAS', '1911', 33125, 828116)
    ,@(17, '1101441905', 'ARNALDO ANDRES BUELVAS MERCADO', '1911', 33125, 828116)
    ,@(18, '1101441905', 'ARNALDO ANDRES BUELVAS MERCADO', '1910', 33125, 828116)
    ,@(19, '1101441905', 'ARNALDO ANDRES BUELVAS MERCADO', '1909', 33125, 828116)
    ,@(20, '9091426', 'LUIS ALBERTO TORRES CABARCAS', '2105', 60000, 1200000)
    ,@(21, '9091426', 'LUIS ALBERTO TORRES CABARCAS', '2104', 60000, 1200000)
    ,@(22, '9091426', 'LUIS ALBERTO TORRES CABARCAS', '2103', 60000, 1200000)
    ,@(23, '9091426', 'LUIS ALBERTO TORRES CABARCAS', '2102', 60000, 1200000)
    ,@(24, '9091426', 'LUIS ALBERTO TORRES CABARCAS', '2101', 60000, 1200000)
    ,@(25, '9091426', 'LUIS ALBERTO TORRES CABARCAS', '2012', 60000, 1200000)
    ,@(26, '1072522889', 'ADRIANA PAOLA PINEDO BELTRAN', '2205', 36341, 2200000)
    ,@(27, '1072522889', 'ADRIANA PAOLA PINEDO BELTRAN', '2204', 36341, 2200000)
    ,@(28, '1072522889', 'ADRIANA PAOLA PINEDO BELTRAN', '2203', 36341, 2200000)
    ,@(29, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2105', 118720, 2968000)
    ,@(30, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2104', 118720, 2968000)
    ,@(31, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2103', 118720, 2968000)
    ,@(32, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2102', 118720, 2968000)
    ,@(33, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2101', 118720, 2968000)
    ,@(34, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2012', 118720, 2968000)
    ,@(35, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2011', 118720, 2968000)
    ,@(36, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2010', 118720, 2968000)
    ,@(37, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2009', 118720, 2968000)
    ,@(38, '1143325816', 'FABIO ENRIQUE PRASCA HERNANDEZ', '2008', 118720, 2968000)
    ,@(39, '92641989', 'JOSE GABRIEL ESPITIA SIERRA', '2001', 87061, 2968000)
    ,@(40, '1050969832', 'YEAN DAVID PABON RIVERA', '2105', 36341, 3260870)
    ,@(41, '1050969832', 'YEAN DAVID PABON RIVERA', '2104', 36341, 3260870)
    ,@(42, '1050969832', 'YEAN DAVID PABON RIVERA', '2103', 36341, 3260870)
    ,@(43, '1050969832', 'YEAN DAVID PABON RIVERA', '2102', 36341, 3260870)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

